# preparation publication 0.2.0
# - bump Version metadata value 0.1.1 -> 0.2.0
# - bump Date metadata value to the new publication timestamp
# - add a new "Jurisdiction" / "iso:code:3166:FR" row right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version (row 3, column B)
$ws.Range("B3").Value = "0.2.0"

# Date (row 8, column B)
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# Insert a new row right after the "Contact" row (row 10) so it becomes row 11,
# pushing Description/Purpose/... etc. down by one row.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Match the formatting used by the rest of the table (inherit from the row above)
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
